# Apply the "filled_template" fixture correction: the header/label row (row 3)
# had field-name placeholders (e.g. "ti2:type") typed into data cells instead of
# valid MODS controlled-vocabulary values. Replace those with correct values,
# tidy the row-3 explicit height override, and reset the saved view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("P3").Value = "abbreviated"
$ws.Range("AF3").Value = "abbreviated"
$ws.Range("AW3").Value = "personal"
$ws.Range("AX3").Value = "primary"
$ws.Range("BW3").Value = "personal"
$ws.Range("CY3").Value = "personal"
$ws.Range("EA3").Value = "personal"
$ws.Range("FC3").Value = "personal"
$ws.Range("GD3").Value = "text"
$ws.Range("GE3").Value = "yes"
$ws.Range("GF3").Value = "text"
$ws.Range("GH3").Value = "text"
$ws.Range("HB3").Value = "yes"
$ws.Range("HC3").Value = "w3cdtf"
$ws.Range("HD3").Value = "approximate"
$ws.Range("HE3").Value = "start"
$ws.Range("HG3").Value = "approximate"
$ws.Range("HH3").Value = "start"
$ws.Range("HJ3").Value = "yes"
$ws.Range("HK3").Value = "w3cdtf"
$ws.Range("HL3").Value = "approximate"
$ws.Range("HM3").Value = "start"
$ws.Range("HO3").Value = "approximate"
$ws.Range("HP3").Value = "start"
$ws.Range("HR3").Value = "yes"
$ws.Range("HS3").Value = "w3cdtf"
$ws.Range("HT3").Value = "approximate"
$ws.Range("HU3").Value = "start"
$ws.Range("HW3").Value = "approximate"
$ws.Range("HX3").Value = "start"
$ws.Range("HZ3").Value = "yes"
$ws.Range("IA3").Value = "w3cdtf"
$ws.Range("IB3").Value = "approximate"
$ws.Range("IC3").Value = "start"
$ws.Range("IE3").Value = "approximate"
$ws.Range("IF3").Value = "start"
$ws.Range("IH3").Value = "yes"
$ws.Range("IJ3").Value = "w3cdtf"
$ws.Range("IK3").Value = "approximate"
$ws.Range("IL3").Value = "start"
$ws.Range("IP3").Value = "approximate"
$ws.Range("IQ3").Value = "start"
$ws.Range("JA3").Value = "marcgac"
$ws.Range("JE3").Value = "continuing"
$ws.Range("JG3").Value = "yes"
$ws.Range("JH3").Value = "w3cdtf"
$ws.Range("JI3").Value = "approximate"
$ws.Range("JJ3").Value = "start"
$ws.Range("JL3").Value = "approximate"
$ws.Range("JM3").Value = "start"
$ws.Range("JO3").Value = "yes"
$ws.Range("JP3").Value = "w3cdtf"
$ws.Range("JQ3").Value = "approximate"
$ws.Range("JR3").Value = "start"
$ws.Range("JT3").Value = "approximate"
$ws.Range("JU3").Value = "start"
$ws.Range("JW3").Value = "yes"
$ws.Range("JX3").Value = "w3cdtf"
$ws.Range("JY3").Value = "approximate"
$ws.Range("JZ3").Value = "start"
$ws.Range("KB3").Value = "approximate"
$ws.Range("KC3").Value = "start"
$ws.Range("KE3").Value = "yes"
$ws.Range("KF3").Value = "w3cdtf"
$ws.Range("KG3").Value = "approximate"
$ws.Range("KH3").Value = "start"
$ws.Range("KJ3").Value = "approximate"
$ws.Range("KK3").Value = "start"
$ws.Range("KM3").Value = "yes"
$ws.Range("KO3").Value = "w3cdtf"
$ws.Range("KP3").Value = "approximate"
$ws.Range("KQ3").Value = "start"
$ws.Range("KU3").Value = "approximate"
$ws.Range("KV3").Value = "start"
$ws.Range("LF3").Value = "marcgac"
$ws.Range("LJ3").Value = "continuing"
$ws.Range("LN3").Value = "rfc3066"
$ws.Range("LS3").Value = "rfc3066"
$ws.Range("LY3").Value = "rfc3066"
$ws.Range("MH3").Value = "access"
$ws.Range("MI3").Value = "born digital"
$ws.Range("ON3").Value = "access"
$ws.Range("OO3").Value = "born digital"
$ws.Range("PA3").Value = "access"
$ws.Range("PB3").Value = "born digital"
$ws.Range("RK3").Value = "personal"
$ws.Range("RS3").Value = "abbreviated"
$ws.Range("TB3").Value = "personal"
$ws.Range("TJ3").Value = "abbreviated"
$ws.Range("UT3").Value = "personal"
$ws.Range("VB3").Value = "abbreviated"
$ws.Range("WL3").Value = "personal"
$ws.Range("WT3").Value = "abbreviated"
$ws.Range("YD3").Value = "personal"
$ws.Range("YL3").Value = "abbreviated"
$ws.Range("AHD3").Value = "preceding"
$ws.Range("AHH3").Value = "personal"
$ws.Range("AHY3").Value = "yes"
$ws.Range("AHZ3").Value = "w3cdtf"
$ws.Range("AIA3").Value = "approximate"
$ws.Range("AIB3").Value = "start"
$ws.Range("AID3").Value = "approximate"
$ws.Range("AIE3").Value = "start"
$ws.Range("AIG3").Value = "yes"
$ws.Range("AIH3").Value = "w3cdtf"
$ws.Range("AII3").Value = "approximate"
$ws.Range("AIJ3").Value = "start"
$ws.Range("AIL3").Value = "approximate"
$ws.Range("AIM3").Value = "start"
$ws.Range("AIQ3").Value = "marcgac"
$ws.Range("AJG3").Value = "preceding"
$ws.Range("AJK3").Value = "personal"
$ws.Range("AKB3").Value = "yes"
$ws.Range("AKC3").Value = "w3cdtf"
$ws.Range("AKD3").Value = "approximate"
$ws.Range("AKE3").Value = "start"
$ws.Range("AKG3").Value = "approximate"
$ws.Range("AKH3").Value = "start"
$ws.Range("AKJ3").Value = "yes"
$ws.Range("AKK3").Value = "w3cdtf"
$ws.Range("AKL3").Value = "approximate"
$ws.Range("AKM3").Value = "start"
$ws.Range("AKO3").Value = "approximate"
$ws.Range("AKP3").Value = "start"
$ws.Range("AKT3").Value = "marcgac"
$ws.Range("ALK3").Value = "preceding"
$ws.Range("ALO3").Value = "personal"
$ws.Range("AMF3").Value = "yes"
$ws.Range("AMG3").Value = "w3cdtf"
$ws.Range("AMH3").Value = "approximate"
$ws.Range("AMI3").Value = "start"
$ws.Range("AMK3").Value = "approximate"
$ws.Range("AML3").Value = "start"
$ws.Range("AMN3").Value = "yes"
$ws.Range("AMO3").Value = "w3cdtf"
$ws.Range("AMP3").Value = "approximate"
$ws.Range("AMQ3").Value = "start"
$ws.Range("AMS3").Value = "approximate"
$ws.Range("AMT3").Value = "start"
$ws.Range("AMX3").Value = "marcgac"
$ws.Range("ANX3").Value = "primary"
$ws.Range("ANZ3").Value = "rfc3066"
$ws.Range("AOJ3").Value = "rfc3066"

# Row 3 no longer needs its explicit ht="60" override - let it size naturally.
$ws.Rows.Item(3).EntireRow.AutoFit()

# Reset the saved sheet view back to the top-left, clearing the stored
# topLeftCell/selection scroll state left over from editing far out at column AOM.
$ws.Range("A1").Select() | Out-Null
